$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (re-ranked MAA accuracy figures + refreshed timestamp) ---
$ws.Range("D2").Value = "maa://24702 (94.48), maa://25390 (96.17), maa://36681 (87.34)"
$ws.Range("L2").Value = "*maa://24633 (56.52), *maa://30515 (69.9), *maa://34787 (73.33), maa://39402 (91.67), ***maa://20792 (11.93), ***maa://29083 (27.78)"
$ws.Range("AF2").Value = "maa://25251 (92.24), ***maa://21730 (25.33), ***maa://39501 (17.24), **maa://36675 (50.0)"
$ws.Range("T3").Value = "maa://24617 (89.74), **maa://20790 (43.48), ***maa://37170 (16.92), maa://45854 (93.75)"
$ws.Range("AF3").Value = "*maa://21289 (72.0)"
$ws.Range("L6").Value = "maa://24839 (99.0)"
$ws.Range("A8").Value = "更新日期：2025.02.15 13:17:50"
$ws.Range("AB9").Value = "maa://28711 (87.39), ***maa://22740 (5.66), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (96.0), *maa://45044 (66.67)"
$ws.Range("D10").Value = "***maa://25695 (18.72), ***maa://34206 (20.0), ***maa://39951 (15.69), ***maa://39243 (28.57), *maa://45271 (57.14)"
$ws.Range("X10").Value = "maa://22301 (97.75), maa://45828 (88.89), maa://22726 (100.0)"
$ws.Range("AF10").Value = "*maa://25021 (54.35), *maa://22733 (60.0), **maa://22761 (50.0)"
$ws.Range("T11").Value = "maa://22747 (93.08), maa://22501 (97.67), *maa://45521 (78.57)"
$ws.Range("D13").Value = "maa://24999 (92.05), maa://36673 (93.24), maa://25001 (85.71)"
$ws.Range("AF13").Value = "**maa://22737 (33.33), maa://39883 (91.43), *maa://39885 (53.33)"
$ws.Range("L14").Value = "maa://26245 (96.69), maa://21288 (96.3), maa://39841 (96.0), maa://36682 (97.44)"
$ws.Range("AB19").Value = "*maa://30709 (65.36), *maa://36668 (57.5)"
$ws.Range("D20").Value = "maa://21432 (90.48), maa://25198 (93.58), *maa://20795 (51.16), maa://36680 (93.94)"
$ws.Range("AF21").Value = "maa://22524 (94.5), *maa://22432 (76.71)"
$ws.Range("H22").Value = "maa://25236 (96.67), **maa://21678 (48.94), **maa://22735 (42.86)"
$ws.Range("X24").Value = "maa://29988 (84.92), maa://23504 (93.1), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (77.63), ***maa://22815 (23.08)"
$ws.Range("AB26").Value = "maa://42235 (94.79)"
$ws.Range("T28").Value = "maa://23263 (95.28), *maa://29765 (63.41)"
$ws.Range("X28").Value = "maa://39929 (90.55), maa://41749 (90.48), ***maa://39723 (13.89)"
$ws.Range("L29").Value = "maa://28432 (93.43), *maa://28440 (79.63), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("X30").Value = "maa://39477 (90.0)"
$ws.Range("H32").Value = "maa://21895 (97.5), maa://36667 (97.53), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("L37").Value = "maa://45718 (98.33), maa://45789 (100.0)"
$ws.Range("AF38").Value = "maa://36697 (86.19)"
$ws.Range("H39").Value = "maa://36670 (88.89), maa://25199 (84.82), maa://30434 (91.25), ***maa://25036 (16.0), maa://45059 (81.25), *maa://44165 (66.67)"
$ws.Range("P39").Value = "maa://24709 (91.39)"
$ws.Range("T39").Value = "maa://45788 (82.61), maa://45790 (81.82)"
$ws.Range("H44").Value = "maa://29768 (98.0), maa://27728 (96.08)"
$ws.Range("H45").Value = "maa://21229 (84.74), maa://30807 (95.65), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (84.21)"
$ws.Range("H53").Value = "maa://32534 (93.92), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (92.28)"
$ws.Range("H62").Value = "maa://42981 (95.0), maa://43903 (100.0)"

# --- New operator row inserted at row 34 (V34:Y34), matching style of neighboring U34 ---
$ws.Range("U34").Copy()
$ws.Range("V34:Y34").PasteSpecial(-4122)
$ws.Range("V34").Value = "诺威尔"
$ws.Range("W34").Value = "-"
$ws.Range("X34").Value = "-"

# --- New trailing row 74, matching style of row 73 ---
$ws.Range("F73:I73").Copy()
$ws.Range("F74:I74").PasteSpecial(-4122)
$ws.Range("F74").Value = "隐德来希"
$ws.Range("G74").Value = "-"
$ws.Range("H74").Value = "-"
